# Regenerate the "K" column (column G) values for rows 2-37 using the
# newly-computed strike values (K instead of Strike#).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 9
    3  = 8
    4  = 6
    5  = 7
    6  = 5
    7  = 5
    8  = 5
    9  = 7
    10 = 9
    11 = 12
    12 = 1
    13 = 12
    14 = 3
    15 = 11
    16 = 2
    17 = 9
    18 = 6
    19 = 5
    20 = 3
    21 = 7
    22 = 7
    23 = 5
    24 = 8
    25 = 10
    26 = 4
    27 = 9
    28 = 6
    29 = 10
    30 = 5
    31 = 7
    32 = 3
    33 = 7
    34 = 4
    35 = 4
    36 = 5
    37 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
